$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-7 (years 2000, 2005, 2006, 2007, 2008, 2009).
# The remaining rows (old 8-11, for years 2010-2013) shift up to become rows 2-5.
$ws.Range("A2:J7").EntireRow.Delete()
